$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix B68: was stored as text "2", should become a real number 2
$ws.Range("B68").Value = 2

# Append new row 69 with the new annotation entry
$ws.Range("A69").Value = "Ruilin"

# B69 keeps "3" stored as text (mirrors how the source data had it),
# so force text formatting before assigning, then drop the format again
# so no extra style sticks to the cell.
$ws.Range("B69").NumberFormat = "@"
$ws.Range("B69").Value = "3"
$ws.Range("B69").ClearFormats()

$ws.Range("C69").Value = "无"
$ws.Range("D69").Value = "DFT"
$ws.Range("E69").Value = "EXP"
$ws.Range("F69").Value = "6325282a-75f6-4567-8bb3-3102657c705c"
$ws.Range("G69").Value = "fm5jfAwPbOfP6_annotated.xlsx"
$ws.Range("H69").Value = "I found the empirical evaluation to be weak."
